# Apply cell value updates from the cryptos list refresh (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.606.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.228.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.46%  "

$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.58"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +17.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0963"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.555.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.856"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.233.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.559.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0965"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "

$ws.Range("E23").Value = "  +15.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.91%  "

$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("E26").Value = "  +2.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.75%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("E32").Value = "  +0.39%  "

$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("E34").Value = "  +2.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0720"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +18.37%  "

$ws.Range("E37").Value = "  -2.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0285"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.64%  "

$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +16.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.00%  "

$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.210"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.12%  "

$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("E50").Value = "  +6.36%  "

$ws.Range("E51").Value = "  +1.36%  "
